# Helper: convert an "RRGGBB" hex string into the BGR-packed integer that
# the ColorFormat.RGB COM property expects (Windows OLE_COLOR / COLORREF
# order is 0x00BBGGRR, i.e. reversed from the familiar RRGGBB hex string).
function HexToOleColor([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation

# --- 1. Re-colour the slide-master theme (currently "Integral" / "Red
#        Violet") so it uses the stock "Office Theme" colour palette. The
#        font scheme and format scheme are already identical between the
#        two themes, so only the 12 colour-scheme slots need to change.
#        ColorScheme.Colors() uses the fixed VBA/COM order:
#          1 dk1  2 lt1  3 dk2  4 lt2  5 accent1 6 accent2
#          7 accent3 8 accent4 9 accent5 10 accent6 11 hlink 12 folHlink
$officeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$sm = $p.SlideMaster
$cs = $sm.ColorScheme
for ($i = 1; $i -le 12; $i++) {
    $cs.Colors($i).RGB = HexToOleColor($officeColors[$i - 1])
}

# --- 2. Re-colour the notes-master theme (currently "Office Theme") so it
#        uses the "Integral" / "Red Violet" palette that used to live on the
#        slide master -- the two themes effectively trade places.
$redVioletColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "454551",  # dk2
    "D8D9DC",  # lt2
    "E32D91",  # accent1
    "C830CC",  # accent2
    "4EA6DC",  # accent3
    "4775E7",  # accent4
    "8971E1",  # accent5
    "D54773",  # accent6
    "6B9F25",  # hlink
    "8C8C8C"   # folHlink
)

$nm = $p.NotesMaster
$ncs = $nm.ColorScheme
for ($i = 1; $i -le 12; $i++) {
    $ncs.Colors($i).RGB = HexToOleColor($redVioletColors[$i - 1])
}

# --- 3. Re-apply the table style on the three tables that were still using
#        the legacy "Table_0" custom style id, switching them to the new
#        built-in style id.
for ($slideIdx = 14; $slideIdx -le 16; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle("{6A6CF36E-EF7D-430F-B3DC-1E3AB5A37477}")
        }
    }
}
